# Update email addresses (shared-string content) across the three sheets that
# reference "demo api" account/contact data. Matches commit "api excel data input".
$wb = $excel.ActiveWorkbook

$wsContAdmin = $wb.Worksheets.Item("Cont adminstrator")
$wsContAdmin.Range("C15").Value = "littlepiglesswt813@automation.33mail.com"

$wsReceptie = $wb.Worksheets.Item("Receptie")
$wsReceptie.Range("B3").Value = "moraritza2714@staffcalendis.33mail.com"
$wsReceptie.Range("B4").Value = "rimmelplus3315@staffcalendis.33mail.com"

$wsAngajati = $wb.Worksheets.Item("Angajati")
$wsAngajati.Range("B2").Value = "marilenajohhjss1523@staffcalendis.33mail.com"
$wsAngajati.Range("B3").Value = "ideaforkih3973@staffcalendis.33mail.com"
$wsAngajati.Range("B4").Value = "boomsie4s2863@staffcalendis.33mail.com"
$wsAngajati.Range("B5").Value = "ocarinass3053@staffcalendis.33mail.com"

# Re-assert the (duplicated, hidden) autofilter-derived defined names that the
# authoring pipeline appends on every regeneration for each autofiltered sheet
# (Cont adminstrator / Domenii / Domenii existente).
$wsDomenii = $wb.Worksheets.Item("Domenii")
$wsDomeniiExistente = $wb.Worksheets.Item("Domenii existente")

$wsContAdmin.Names.Add("_xlnm._FilterDatabase_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0", "='Cont adminstrator'!$A$1:$A$19")
$wsDomenii.Names.Add("_xlnm._FilterDatabase_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0", "=Domenii!$A$4:$A$7")
$wsDomeniiExistente.Names.Add("_xlnm._FilterDatabase_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0", "='Domenii existente'!$A$1:$Q$15")

# Move the active workbook tab from "Receptie" to "Permisiuni" (activeTab 1 -> 7).
$wsPermisiuni = $wb.Worksheets.Item("Permisiuni")
$wsPermisiuni.Activate()

